# Update crypto price/volume data as per the daily GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.955.13"
$ws.Range("E2").Value = "'  +2.54%  "
$ws.Range("D3").Value = "'2.227.05"
$ws.Range("E3").Value = "'  +0.91%  "
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'294.47"
$ws.Range("E5").Value = "'  -0.14%  "
$ws.Range("D6").Value = "'85.97"
$ws.Range("E6").Value = "'  +7.67%  "
$ws.Range("E7").Value = "'  +2.09%  "
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("D9").Value = "'0.471"
$ws.Range("E9").Value = "'  +3.53%  "
$ws.Range("D10").Value = "'31.02"
$ws.Range("E10").Value = "'  +11.59%  "
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "'  +2.60%  "
$ws.Range("D12").Value = "'46.89"
$ws.Range("E12").Value = "'  +2.05%  "
$ws.Range("E13").Value = "'  +1.30%  "
$ws.Range("E14").Value = "'  +6.04%  "
$ws.Range("D15").Value = "'2.575.83"
$ws.Range("E15").Value = "'  +1.06%  "
$ws.Range("B16").Value = "'WrappedEther"
$ws.Range("C16").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'2.340.30"
$ws.Range("E16").Value = "'  +5.39%  "
$ws.Range("B17").Value = "'Chainlink"
$ws.Range("C17").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'14.04"
$ws.Range("E17").Value = "'  +1.87%  "
$ws.Range("D18").Value = "'0.727"
$ws.Range("E18").Value = "'  +3.17%  "
$ws.Range("D19").Value = "'39.866.69"
$ws.Range("E19").Value = "'  +2.57%  "
$ws.Range("E20").Value = "'  +3.90%  "
$ws.Range("D21").Value = "'5.78"
$ws.Range("E21").Value = "'  +2.49%  "
$ws.Range("D22").Value = "'10.77"
$ws.Range("E22").Value = "'  +10.45%  "
$ws.Range("D23").Value = "'65.13"
$ws.Range("E23").Value = "'  +1.08%  "
$ws.Range("D24").Value = "'234.96"
$ws.Range("E24").Value = "'  +5.30%  "
$ws.Range("E25").Value = "'  -0.14%  "
$ws.Range("D26").Value = "'2.46"
$ws.Range("E26").Value = "'  +3.59%  "
$ws.Range("E27").Value = "'  +5.83%  "
$ws.Range("D28").Value = "'22.73"
$ws.Range("E28").Value = "'  +2.60%  "
$ws.Range("E29").Value = "'  +2.98%  "
$ws.Range("D30").Value = "'9.20"
$ws.Range("E30").Value = "'  +4.16%  "
$ws.Range("D31").Value = "'33.09"
$ws.Range("E31").Value = "'  +7.01%  "
$ws.Range("D32").Value = "'151.99"
$ws.Range("E32").Value = "'  +2.90%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "'  -0.03%  "
$ws.Range("E34").Value = "'  +2.98%  "
$ws.Range("D35").Value = "'0.0716"
$ws.Range("E35").Value = "'  +5.19%  "
$ws.Range("E36").Value = "'  +2.39%  "
$ws.Range("D37").Value = "'16.22"
$ws.Range("E37").Value = "'  +14.30%  "
$ws.Range("E38").Value = "'  +2.72%  "
$ws.Range("D39").Value = "'0.0996"
$ws.Range("E39").Value = "'  +3.88%  "
$ws.Range("D40").Value = "'2.71"
$ws.Range("E40").Value = "'  +4.02%  "
$ws.Range("E41").Value = "'  +6.84%  "
$ws.Range("D42").Value = "'3.79"
$ws.Range("E42").Value = "'  +6.49%  "
$ws.Range("D43").Value = "'2.028.44"
$ws.Range("E43").Value = "'  +7.36%  "
$ws.Range("D44").Value = "'2.20"
$ws.Range("E44").Value = "'  +7.01%  "
$ws.Range("D45").Value = "'0.0269"
$ws.Range("E45").Value = "'  +6.43%  "
$ws.Range("D46").Value = "'9.90"
$ws.Range("E46").Value = "'  +13.13%  "
$ws.Range("D47").Value = "'16.18"
$ws.Range("E47").Value = "'  +2.01%  "
$ws.Range("E48").Value = "'  +2.84%  "
$ws.Range("D49").Value = "'2.449.10"
$ws.Range("E49").Value = "'  +1.18%  "
$ws.Range("D50").Value = "'70.61"
$ws.Range("E50").Value = "'  +1.23%  "
$ws.Range("E51").Value = "'  +14.92%  "
